$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column: header in H1, matching the formatting of the
# neighboring header cells (e.g. G1) by copying its format onto H1.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column data for rows 2-4
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
